# "база данных.xlsx" — rebuild the "расписание" sheet as a Mon-Fri grid.
# The original sheet listed one row per lesson, grouped under a day name in
# column B every 7 rows (понедельник/вторник/среда/четверг/пятница).
# The new layout pivots this into a single table: column headers are the
# five week days (row 1), and each subsequent row is one time slot shared
# across all days (row 2 = 09:05-10:25, ... row 8 = 15:20-16:xx).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# ---- 1. Formats -----------------------------------------------------
# Header row (A1:B1) is already bold/centered (style used for "расписание"
# / "ключевое слово"); stamp the same format across C1:F1 so every day
# header (and the trailing blank F1) matches it.
$ws.Range("B1").Copy()
$ws.Range("C1:F1").PasteSpecial(-4122)

# Body cells (A2:E8) take the plain "главный" look already used by A2.
$ws.Range("A2").Copy()
$ws.Range("A2:E8").PasteSpecial(-4122)

# Two cells (the Monday 12:30 slot and the Monday 15:20 slot) keep the
# alternate font that "B9" (a day-name cell that will be cleared below)
# still carries at this point in the script.
$ws.Range("B9").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A8").PasteSpecial(-4122)

# ---- 2. Values --------------------------------------------------------
$ws.Range("A1").Value = "расписание на понедельник:"
$ws.Range("B1").Value = "расписание на вторник:"
$ws.Range("C1").Value = "расписание на среду:"
$ws.Range("D1").Value = "расписание на четверг:"
$ws.Range("E1").Value = "расписание на пятницу:"
$ws.Range("F1").Value = ""

$grid = @(
    @("09:05-10:25   главный", "09:05-10:25   главный", "09:05-10:25   главный", "09:05-10:25   главный", "09:05-10:25   главный"),
    @("10:40-11:25   обществознание", "10:40-11:25   математика", "10:40-11:25   химия", "10:40-11:25   биология", "10:40-11:25   математика"),
    @("11:40-12:20   обществознание", "11:40-12:20   химия", "11:40-12:20   экономика", "11:40-12:20   английский", "11:40-12:20   математика"),
    @("12:30-13:10   русский", "12:30-13:10   история", "12:30-13:10   английский", "12:30-13:10   эвритмия", "12:30-13:10   история"),
    @("13:30-14:10   живопись", "13:30-14:10   информатика", "13:30-14:10   информатика", "13:30-14:10   математика", "13:30-14:10   экономика"),
    @("14:30-15:10   искусство", "14:30-15:10   английский", "14:30-15:10   английский", "14:30-15:10   проект", "14:30-15:10   английский"),
    @("15:20-16:10   репетиция", "15:20-16:10   искусство", "15:20-16:00   физкультура", "15:20-16:00   проект", "15:20-16:00   русский")
)

for ($r = 0; $r -lt $grid.Length; $r++) {
    $rowVals = $grid[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item(2 + $r, 1 + $c).Value = $rowVals[$c]
    }
}

# ---- 3. Drop the old per-day rows -------------------------------------
# Rows 9-36 held the remaining four days (вторник..пятница) one lesson per
# row; that data now lives in columns B:E of rows 2-8, so clear it. Rows
# without an explicit height (9-18) disappear entirely once emptied; rows
# 19-36 already carry ht="15.75" and stay behind as blank stub rows.
$ws.Range("A9:E36").Clear()

# ---- 4. Column widths ---------------------------------------------
# Columns A:E widen to fit the day tables; F onward keeps the default.
for ($c = 1; $c -le 5; $c++) {
    $ws.Columns.Item($c).ColumnWidth = 26.857142857142858
}
